$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $findText"
    }
    return $ok
}

# 1. Merge paragraphs "So what is summarization..." + "There are two main strategies..."
Replace-Text "overall meaning. ^pThere are two main strategies" "overall meaning. There are two main strategies"

# 1b. Merge with "The goal of this project..." paragraph, inserting ". "
Replace-Text "for the summary^pThe goal of this project" "for the summary. The goal of this project"

# 2a. "Text summarization is becoming" -> "It is becoming"
Replace-Text "-Text summarization is becoming" "-It is becoming"

# 2b. Merge "business analysis." paragraph with "With this technique..." paragraph, adding a space
Replace-Text "business analysis.^pWith this technique" "business analysis. With this technique"

# 2c. Delete the "Here is a scenario..." paragraph entirely (keep the following empty paragraph)
Replace-Text "every word.^pHere is a scenario that your boss wants you to write a summary of the report to save his time. Instead of reading the whole report word by word, you can use summarization model to generate a summary and give that to your boss. It is quite convenient and efficient.^p" "every word.^p"

# 3a. "multiple summaries. This data is shown in the table." -> "multiple summaries and this data is shown in the table."
Replace-Text "multiple summaries. This data is shown in the table." "multiple summaries and this data is shown in the table."

# 3b. Delete "The articles have an average..." paragraph entirely, merging into previous paragraph
Replace-Text "in the table.^pThe articles have an average of 781 tokens while the summaries have an average of 56 tokens." "in the table."

# 3c. Merge with "For pre-processing step..." paragraph
Replace-Text "in the table.^pFor pre-processing step" "in the table. For pre-processing step"

# 3d. "noisy words that are not related to news articles. We also lower" -> "noisy words, and also lower"
Replace-Text "noisy words that are not related to news articles. We also lower" "noisy words, and also lower"

# 3e. Merge "...in all words." with "And to save computation time..." paragraph
Replace-Text "in all words.^pAnd to save computation time" "in all words. To save computation time"

# 4. Remove clause about decoder/target vocabulary
Replace-Text "Neural Machine Translation and summaries are generated from the decoder, using target vocabulary. This model" "Neural Machine Translation. This model"

# 5a. Remove "a cleaned version of " before Common Crawl
Replace-Text "C4 dataset which is a cleaned version of Common Crawl" "C4 dataset which is Common Crawl"

# 5b. Remove sentence about T5 Masked Language Model
Replace-Text "text strings. T5 uses the same Masked Language Model as BERT but it is different from BERT-based models that can only output either a class label or a span of the input. The advantage" "text strings. The advantage"

# 6a. Remove clause "which is commonly used in summarization tasks" and merge paragraphs
Replace-Text "F1 from the ROUGE metric which is commonly used in summarization tasks. ^pWe select only" "F1 from the ROUGE metric. We select only"

# 6b. Delete "In this example, we can see..." paragraph and the following empty paragraph entirely
Replace-Text "between the two summaries. ^pIn this example, we can see that there are total 7 words in the system summary, and 6 of them are overlapping words so the precision score should be 6 divided by 7 and recall score is 6 divided by 6. And these two numbers will be used for calculating F1 score for ROUGE-1^p^p-Here is our result." "between the two summaries. ^p-Here is our result."

# 7. Merge "...compared with other models." with "The performance of Transformer-based..." (no space added)
Replace-Text "compared with other models.^pThe performance of Transformer-based" "compared with other models.The performance of Transformer-based"

# 8a. "The accuracy of text summarization models can be increased by" -> "the performance of text summarization models can be improved by"
Replace-Text "The accuracy of text summarization models can be increased by" "the performance of text summarization models can be improved by"

# 8b. Merge "we conclude that " with "the performance..." paragraph
Replace-Text "we conclude that ^pthe performance" "we conclude that the performance"

# 8c. Merge "...transformer model." with "Our current models..." paragraph
Replace-Text "transformer model.^pOur current models" "transformer model. Our current models"

# 9. Merge "...challenging topic." with "The overall accuracy of our models is..." and rephrase
Replace-Text "is still a challenging topic.^pThe overall accuracy of our models is" "is still a challenging topic as the overall accuracy of our models is"

Write-Host "Done"
